# Add a new DNA-repeat-sequence problem row to the "哈希" (Hash) worksheet,
# matching the commit "repeat DNA design  with hash".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Clone row 2 formatting (font/alignment/wrap) down into the new rows 3-8
# so the new cells pick up the same cell style used for existing data rows.
$ws.Range("A2:G2").Copy()
$ws.Range("A3:G8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 3: new problem entry (No. 2, leetcode 187 - "Repeated DNA Sequences")
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 187
$ws.Range("C3").Value = "所有 DNA 都由一系列缩写为 A，C，G 和 T 的核苷酸组成，例如：“ACGAATTCCG”。在研究 DNA 时，识别 DNA 中的重复序列有时会对研究非常有帮助"  # 题目 (problem statement)
$ws.Range("D3").Value = "1 窗口长度是10`n2 从原字符串0索引开始滑动窗口，落入窗口的子字符串是否在迭代集合中`n3 存在，说明这个子字符串出现了至少一次了，放入结果集合中`n4 不存在，说明这个子字符串只出现了一次，就将其放入迭代集合中`n4 迭代完整个原字符串【注意窗口位置】`n5 返回结果字符串，注意去重"  # 解题方法 (approach)
$ws.Range("E3").Value = "集合去重`n滑动窗口"  # 解题关键词 (keywords)
$ws.Range("F3").Value = "窗口长度L，原字符串长度N，N-L+1个子字符串，每个子字符串长度是L。`n时间复杂度：O((N-L+1)L),题目已假定L=10，即时间复杂度是O(N)        "  # 时间复杂度 (time complexity)
$ws.Range("G3").Value = "空间复杂度：集合保存这些子字符串，需要的空间是(N-L+1)L。`nO((N-L+1)L)，题目已假定L=10，即空间复杂度是O(N)"  # 空间复杂度 (space complexity)

# Row heights: row 3 holds the large wrapped paragraphs, rows 4-8 are the
# short blank placeholder rows that were left under it.
$ws.Rows.Item(3).RowHeight = 220
$ws.Range("A4:G8").RowHeight = 21

# Move the view/selection the way the author left it: scrolled to row 2,
# with F7 selected.
$ws.Activate()
$ws.Range("F7").Select()
try { $excel.ActiveWindow.ScrollRow = 2 } catch {}

Write-Output "Added DNA repeat-sequence row to 哈希 sheet"
